$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp update (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 17:05"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1531737
$ws.Range("C4").Value = 4073
$ws.Range("D4").Value = 346786
$ws.Range("E4").Value = 1093890
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 91061

# --- Singapur (row 30) ---
$ws.Range("D30").Value = 9835
$ws.Range("E30").Value = 18486

# --- Filipinas / Republica Dominicana swap places (rows 44-45) ---
# Row 44 now holds "Republica Dominicana" with updated stats
$ws.Range("A44").Value = "Republica Dominicana"
$ws.Range("B44").Value = 12725
$ws.Range("C44").Value = 411
$ws.Range("D44").Value = 6613
$ws.Range("E44").Value = 5678
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 6
$ws.Range("H44").Value = 434

# Row 45 now holds "Filipinas" (its previous stats, unchanged values)
$ws.Range("A45").Value = "Filipinas"
$ws.Range("B45").Value = 12718
$ws.Range("C45").Value = 205
$ws.Range("D45").Value = 2729
$ws.Range("E45").Value = 9158
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 831

# --- Sri Lanka (row 106) ---
$ws.Range("B106").Value = 986
$ws.Range("C106").Value = 5
$ws.Range("E106").Value = 418

# --- Libano (row 108) ---
$ws.Range("D108").Value = 251
$ws.Range("E108").Value = 654
